$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1513.3529
$ws.Range("I40").Value = 1396.7
$ws.Range("J40").Value = 1680
$ws.Range("K40").Value = 1396.7
$ws.Range("L40").Value = 1680
$ws.Range("M40").Value = -1221.7
$ws.Range("N40").Value = -2030

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1114.9231
$ws.Range("I41").Value = 398.8889
$ws.Range("J41").Value = 2726
$ws.Range("K41").Value = 398.8889
$ws.Range("L41").Value = 2726
$ws.Range("M41").Value = 41.11110000000002
$ws.Range("N41").Value = -3606

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2500.4443
$ws.Range("J62").Value = 1200
$ws.Range("L62").Value = 1200
$ws.Range("N62").Value = -2448

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2500.4443
$ws.Range("J65").Value = 1200
$ws.Range("L65").Value = 6000
$ws.Range("N65").Value = -12240

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4501.143
$ws.Range("I69").Value = 3671
$ws.Range("J69").Value = 5123.75
$ws.Range("K69").Value = 11013
$ws.Range("L69").Value = 15371.25
$ws.Range("M69").Value = -10139
$ws.Range("N69").Value = -17119.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 4501.143
$ws.Range("I72").Value = 3671
$ws.Range("J72").Value = 5123.75
$ws.Range("K72").Value = 33039
$ws.Range("L72").Value = 46113.75
$ws.Range("M72").Value = -28671
$ws.Range("N72").Value = -54849.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1966.1923
$ws.Range("J112").Value = 2094.913
$ws.Range("L112").Value = 6284.739
$ws.Range("N112").Value = -8500.739

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2155.7312
$ws.Range("I138").Value = 1100.5625
$ws.Range("J138").Value = 3281.2444
$ws.Range("K138").Value = 3301.6875
$ws.Range("L138").Value = 9843.733200000001
$ws.Range("M138").Value = 1838.3125
$ws.Range("N138").Value = -20123.7332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15082.683
$ws.Range("I32").Value = 15025.189
$ws.Range("K32").Value = 15025.189
$ws.Range("M32").Value = -14738.189

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 38377
$ws.Range("J109").Value = 38377
$ws.Range("L109").Value = 38377
$ws.Range("N109").Value = -41151

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6117.967
$ws.Range("I132").Value = 7989.6113
$ws.Range("J132").Value = 3310.5
$ws.Range("K132").Value = 23968.8339
$ws.Range("L132").Value = 9931.5
$ws.Range("M132").Value = -21438.8339
$ws.Range("N132").Value = -14991.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 127874.81
$ws.Range("I86").Value = 2768.4614
$ws.Range("J86").Value = 670002.3
$ws.Range("K86").Value = 2768.4614
$ws.Range("L86").Value = 670002.3
$ws.Range("M86").Value = -1645.4614
$ws.Range("N86").Value = -672248.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 127874.81
$ws.Range("I89").Value = 2768.4614
$ws.Range("J89").Value = 670002.3
$ws.Range("K89").Value = 13842.307
$ws.Range("L89").Value = 3350011.5
$ws.Range("M89").Value = -8226.307000000001
$ws.Range("N89").Value = -3361243.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 15061.158
$ws.Range("I107").Value = 18944.207
$ws.Range("K107").Value = 18944.207
$ws.Range("M107").Value = -17024.207

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 221.5
$ws.Range("I7").Value = 176.66667
$ws.Range("J7").Value = 240.71428
$ws.Range("K7").Value = 176.66667
$ws.Range("L7").Value = 240.71428
$ws.Range("M7").Value = -63.66667000000001
$ws.Range("N7").Value = -466.71428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 85500.836
$ws.Range("I62").Value = 102061
$ws.Range("K62").Value = 102061
$ws.Range("M62").Value = -101437

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 85500.836
$ws.Range("I65").Value = 102061
$ws.Range("K65").Value = 510305
$ws.Range("M65").Value = -507185

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 31998.875
$ws.Range("J68").Value = 31998.875
$ws.Range("L68").Value = 31998.875
$ws.Range("N68").Value = -33496.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 31998.875
$ws.Range("J71").Value = 31998.875
$ws.Range("L71").Value = 95996.625
$ws.Range("N71").Value = -103484.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5547.8647
$ws.Range("I3").Value = 2093.1333
$ws.Range("J3").Value = 7903.364
$ws.Range("K3").Value = 6279.3999
$ws.Range("L3").Value = 23710.092
$ws.Range("M3").Value = -6167.3999
$ws.Range("N3").Value = -23934.092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 560.7895
$ws.Range("I34").Value = 280
$ws.Range("J34").Value = 690.38464
$ws.Range("K34").Value = 840
$ws.Range("L34").Value = 2071.15392
$ws.Range("M34").Value = -756
$ws.Range("N34").Value = -2239.15392

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2313.6155
$ws.Range("I39").Value = 800
$ws.Range("J39").Value = 2767.7
$ws.Range("K39").Value = 2400
$ws.Range("L39").Value = 8303.099999999999
$ws.Range("M39").Value = -2106
$ws.Range("N39").Value = -8891.099999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3469.0715
$ws.Range("J55").Value = 3469.0715
$ws.Range("L55").Value = 10407.2145
$ws.Range("N55").Value = -10761.2145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 382.84848
$ws.Range("I107").Value = 437.10526
$ws.Range("J107").Value = 309.2143
$ws.Range("K107").Value = 1311.31578
$ws.Range("L107").Value = 927.6428999999999
$ws.Range("M107").Value = 608.6842200000001
$ws.Range("N107").Value = -4767.6429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13890415
$ws.Range("I131").Value = 381.66666
$ws.Range("J131").Value = 15153145
$ws.Range("K131").Value = 1144.99998
$ws.Range("L131").Value = 45459435
$ws.Range("M131").Value = 3895.00002
$ws.Range("N131").Value = -45469515

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6026.1274
$ws.Range("I70").Value = 5764.1665
$ws.Range("J70").Value = 6883.4546
$ws.Range("K70").Value = 5764.1665
$ws.Range("L70").Value = 6883.4546
$ws.Range("M70").Value = -5494.1665
$ws.Range("N70").Value = -7423.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6026.1274
$ws.Range("I73").Value = 5764.1665
$ws.Range("J73").Value = 6883.4546
$ws.Range("K73").Value = 5764.1665
$ws.Range("L73").Value = 6883.4546
$ws.Range("M73").Value = -4828.1665
$ws.Range("N73").Value = -8755.454600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2072
$ws.Range("I113").Value = 1538.5714
$ws.Range("J113").Value = 2818.8
$ws.Range("K113").Value = 1538.5714
$ws.Range("L113").Value = 2818.8
$ws.Range("M113").Value = 631.4286
$ws.Range("N113").Value = -7158.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2306.4546
$ws.Range("I122").Value = 1736.375
$ws.Range("J122").Value = 3826.6667
$ws.Range("K122").Value = 5209.125
$ws.Range("L122").Value = 11480.0001
$ws.Range("M122").Value = -2759.125
$ws.Range("N122").Value = -16380.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 8619.764999999999
$ws.Range("J123").Value = 8619.764999999999
$ws.Range("L123").Value = 8619.764999999999
$ws.Range("N123").Value = -13519.765

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5164
$ws.Range("I7").Value = 5398.8
$ws.Range("K7").Value = 5398.8
$ws.Range("M7").Value = -5286.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1462.7142
$ws.Range("I22").Value = 1295.1666
$ws.Range("J22").Value = 2468
$ws.Range("K22").Value = 1295.1666
$ws.Range("L22").Value = 2468
$ws.Range("M22").Value = -1000.1666
$ws.Range("N22").Value = -3058

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1462.7142
$ws.Range("I27").Value = 1295.1666
$ws.Range("J27").Value = 2468
$ws.Range("K27").Value = 1295.1666
$ws.Range("L27").Value = 2468
$ws.Range("M27").Value = -1188.1666
$ws.Range("N27").Value = -2682

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 19569752
$ws.Range("I122").Value = 27781616
$ws.Range("J122").Value = 14290696
$ws.Range("K122").Value = 83344848
$ws.Range("L122").Value = 42872088
$ws.Range("M122").Value = -83342398
$ws.Range("N122").Value = -42876988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5164
$ws.Range("I126").Value = 5398.8
$ws.Range("K126").Value = 16196.4
$ws.Range("M126").Value = -13726.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 28859068
$ws.Range("I136").Value = 41668204
$ws.Range("J136").Value = 911863.9399999999
$ws.Range("K136").Value = 125004612
$ws.Range("L136").Value = 2735591.82
$ws.Range("M136").Value = -125002062
$ws.Range("N136").Value = -2740691.82

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4181.909
$ws.Range("I62").Value = 3200.2
$ws.Range("K62").Value = 3200.2
$ws.Range("M62").Value = -2576.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4181.909
$ws.Range("I65").Value = 3200.2
$ws.Range("K65").Value = 16001
$ws.Range("M65").Value = -12881
